$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "mani.dubey" URL as the (soon to be) last row, BEFORE the
# row-2 insert below, so it lands on row 22 once everything shifts down.
# Give it an explicit black font color (matches the new font/style added
# to styles.xml in the target workbook).
$ws.Range("A21").Value = "https://kivihealth.com/iam/mani.dubey.4786"
$ws.Range("A21").Font.Color = 0

# Insert a new row at row 2 and populate it with the "vaibhav.nepalia" URL.
# This pushes all the previously existing data rows down by one.
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "https://kivihealth.com/iam/vaibhav.nepalia.3546"

# Reflect the new selection/active cell state.
$ws.Range("A2").Select() | Out-Null

# Switch the sheet's print orientation to portrait.
$ws.PageSetup.Orientation = 1
